$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Append the new mail-log row (row 8)
$ws.Range("A8").Value = "Sollicitatie marketingfunctie"
$ws.Range("B8").Value = "mailmind.test@zohomail.eu"
$ws.Range("C8").Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$ws.Range("D8").Value = "Overig"
$ws.Range("F8").Value = "2025-06-19 11:28:10"
$ws.Range("G8").Value = "Nee"

# Extend the conditional formatting ranges to cover the new row
$dFc = $ws.Range("D2:D7").FormatConditions
for ($i = 1; $i -le $dFc.Count; $i++) {
    $dFc.Item($i).ModifyAppliesToRange($ws.Range("D2:D8"))
}
$gFc = $ws.Range("G2:G7").FormatConditions
for ($i = 1; $i -le $gFc.Count; $i++) {
    $gFc.Item($i).ModifyAppliesToRange($ws.Range("G2:G8"))
}

# Update the Dashboard summary count for "Overig" (4 -> 5)
$dash.Range("B2").Value = 5
